$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1931818181818182
$ws.Range("C2").Value = 0.5643939393939394
$ws.Range("J2").Value = 0.01136363636363636
$ws.Range("P2").Value = 0.1628787878787879
$ws.Range("S2").Value = 0.06818181818181818
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.07975460122699386
$ws.Range("J3").Value = 0.03067484662576687
$ws.Range("P3").Value = 0.7116564417177914
$ws.Range("S3").Value = 0.1717791411042945
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.08040201005025126
$ws.Range("D6").Value = 0.02010050251256281
$ws.Range("F6").Value = 0.04522613065326633
$ws.Range("J6").Value = 0.2412060301507538
$ws.Range("O6").Value = 0.02512562814070352
$ws.Range("Q6").Value = 0.1809045226130653
$ws.Range("R6").Value = 0.09547738693467336
$ws.Range("S6").Value = 0.3115577889447236
$ws.Range("B7").Value = 0.09583333333333334
$ws.Range("D7").Value = 0.01666666666666667
$ws.Range("E7").Value = 0.004166666666666667
$ws.Range("F7").Value = 0.075
$ws.Range("J7").Value = 0.1375
$ws.Range("O7").Value = 0.02916666666666667
$ws.Range("Q7").Value = 0.1625
$ws.Range("R7").Value = 0.1083333333333333
$ws.Range("S7").Value = 0.3708333333333333
$ws.Range("B8").Value = 0.09544468546637744
$ws.Range("D8").Value = 0.01735357917570499
$ws.Range("F8").Value = 0.06724511930585683
$ws.Range("J8").Value = 0.1106290672451193
$ws.Range("O8").Value = 0.03904555314533623
$ws.Range("Q8").Value = 0.1691973969631236
$ws.Range("R8").Value = 0.09761388286334056
$ws.Range("S8").Value = 0.403470715835141
$ws.Range("B9").Value = 0.08695652173913043
$ws.Range("D9").Value = 0.03623188405797102
$ws.Range("F9").Value = 0.03623188405797102
$ws.Range("J9").Value = 0.108695652173913
$ws.Range("O9").Value = 0.007246376811594203
$ws.Range("Q9").Value = 0.2173913043478261
$ws.Range("R9").Value = 0.08695652173913043
$ws.Range("S9").Value = 0.4202898550724637
$ws.Range("B10").Value = 0.09497206703910614
$ws.Range("D10").Value = 0.02154828411811652
$ws.Range("E10").Value = 0.0007980845969672786
$ws.Range("F10").Value = 0.06464485235434957
$ws.Range("J10").Value = 0.1053471667996808
$ws.Range("O10").Value = 0.01995211492418196
$ws.Range("Q10").Value = 0.2266560255387071
$ws.Range("R10").Value = 0.1077414205905826
$ws.Range("S10").Value = 0.3583399840383081
$ws.Range("G11").Value = 0.1408045977011494
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.1752873563218391
$ws.Range("L11").Value = 0.5804597701149425
$ws.Range("S11").Value = 0.02011494252873563
$ws.Range("G12").Value = 0.775609756097561
$ws.Range("J12").Value = 0.1609756097560976
$ws.Range("K12").Value = 0.01951219512195122
$ws.Range("L12").Value = 0.00975609756097561
$ws.Range("S12").Value = 0.03414634146341464
$ws.Range("G13").Value = 0.6896551724137931
$ws.Range("J13").Value = 0.2413793103448276
$ws.Range("S13").Value = 0.06896551724137931
$ws.Range("F15").Value = 0.01834862385321101
$ws.Range("H15").Value = 0.1559633027522936
$ws.Range("I15").Value = 0.04587155963302753
$ws.Range("J15").Value = 0.3623853211009174
$ws.Range("K15").Value = 0.06880733944954129
$ws.Range("M15").Value = 0.01376146788990826
$ws.Range("O15").Value = 0.02752293577981652
$ws.Range("S15").Value = 0.3073394495412844
$ws.Range("F16").Value = 0.005319148936170213
$ws.Range("H16").Value = 0.1702127659574468
$ws.Range("I16").Value = 0.101063829787234
$ws.Range("J16").Value = 0.3191489361702128
$ws.Range("K16").Value = 0.1436170212765958
$ws.Range("M16").Value = 0.03191489361702127
$ws.Range("O16").Value = 0.09042553191489362
$ws.Range("S16").Value = 0.1382978723404255
$ws.Range("F17").Value = 0.0128755364806867
$ws.Range("H17").Value = 0.1652360515021459
$ws.Range("I17").Value = 0.05793991416309013
$ws.Range("J17").Value = 0.4377682403433477
$ws.Range("K17").Value = 0.1330472103004292
$ws.Range("M17").Value = 0.01502145922746781
$ws.Range("O17").Value = 0.06866952789699571
$ws.Range("S17").Value = 0.1094420600858369
$ws.Range("F18").Value = 0.02521008403361345
$ws.Range("H18").Value = 0.180672268907563
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("J18").Value = 0.3991596638655462
$ws.Range("K18").Value = 0.08823529411764706
$ws.Range("M18").Value = 0.02941176470588235
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("S18").Value = 0.134453781512605
$ws.Range("F19").Value = 0.01229508196721311
$ws.Range("H19").Value = 0.230327868852459
$ws.Range("I19").Value = 0.05491803278688524
$ws.Range("J19").Value = 0.3778688524590164
$ws.Range("K19").Value = 0.1295081967213115
$ws.Range("M19").Value = 0.02868852459016394
$ws.Range("N19").Value = 0.000819672131147541
$ws.Range("O19").Value = 0.05491803278688524
$ws.Range("S19").Value = 0.110655737704918
